$wb = $excel.ActiveWorkbook

# --- Create the new "SPRINT 2" sheet by duplicating "SPRINT 1" so that all
#     formatting (column widths, date styles, borders, page margins, etc.)
#     carries over exactly, then overwrite the content that differs. ---
$sprint1 = $wb.Worksheets.Item("SPRINT 1")
$sprint1.Copy($null, $sprint1)
$ws = $wb.Worksheets.Item($sprint1.Index + 1)
$ws.Name = "SPRINT 2"

# Row 3 - Edubirdie Etat art
$ws.Range("B3").Value = 43923
$ws.Range("C3").Value = "Edubirdie Etat art"
$ws.Range("D3").Value = 1

# Row 4 - Prepostseo Etat art
$ws.Range("B4").Value = 43923
$ws.Range("C4").Value = "Prepostseo Etat art"
$ws.Range("D4").Value = 2

# Row 5 - plagiarismsearch Etat art
$ws.Range("B5").Value = 43926
$ws.Range("C5").Value = "plagiarismsearch Etat art"
$ws.Range("D5").Value = 2

# Row 6 - install XAMPP
$ws.Range("B6").Value = 43926
$ws.Range("C6").Value = "install XAMPP"
$ws.Range("D6").Value = 0.5

# Row 7 - Copyleaks test Etat art (trailing space preserved)
$ws.Range("B7").Value = 43926
$ws.Range("C7").Value = "Copyleaks test Etat art "
$ws.Range("D7").Value = 3.5

# Row 8 - Unicheck Etat art
$ws.Range("B8").Value = 43926
$ws.Range("C8").Value = "Unicheck Etat art"
$ws.Range("D8").Value = 0.5

# Rows 9-11 carried date/text/time data from SPRINT 1 that SPRINT 2 doesn't
# need - clear it, leaving just the styled (empty) date cells in column B.
$ws.Range("B9:D11").ClearContents()

# Row 18 ("Sprint 1 review..." note) isn't part of SPRINT 2 - remove it so
# the row disappears entirely and the used range shrinks back to row 17.
$ws.Range("C18").ClearContents()

$ws.Range("D17").Formula = "=SUM(D3:D11)"

# SPRINT 1 is no longer the active tab; its saved selection becomes a
# select-all (as if the corner "select all" button/Ctrl+A was used there
# before switching away).
$sprint1.Activate()
$sprint1.Cells.Select()

# Leave "SPRINT 2" as the active tab (it's the newly-created, now-last
# sheet) with the saved selection/active cell shown in the sheet.
$ws.Activate()
$ws.Range("C9").Select()
